$wb = $excel.ActiveWorkbook

# Sheet "展览" - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 107
$ws1.Range("F5").Value = 164
$ws1.Range("F6").Value = 9418
$ws1.Range("F7").Value = 847
$ws1.Range("F10").Value = 1149
$ws1.Range("F11").Value = 150
$ws1.Range("F15").Value = 424

# Sheet "全部类型" - update column F (想去人数) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 107
$ws4.Range("F6").Value = 164
$ws4.Range("F7").Value = 9418
$ws4.Range("F8").Value = 847
$ws4.Range("F11").Value = 1149
$ws4.Range("F12").Value = 150
$ws4.Range("F16").Value = 424
